$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(
    -0.22802204078956834,
    -0.14176289238021766,
    -0.07864751493346667,
    -0.07027534696283233,
    -0.06703326481512395,
    -0.04105750762648519,
    -0.03079135764950225,
    -0.020737997446487633,
    -0.018672225145955235,
    -0.01662622425803839,
    -0.013623837900113678,
    0.006701319791543536,
    0.010234929628353484,
    0.018248701119323307,
    0.01926002238498281,
    -0.0060330508073080935,
    -0.004003441776549543,
    -0.08009121387007667,
    -0.07593383733311798,
    -0.008017353860367393,
    -0.0040057871085990016,
    -0.05878483347342289,
    -0.05344285938826321,
    -0.02009815994110742,
    -0.023174161632306678,
    -0.02066183851240666,
    -0.018148847165163406,
    -0.01608731227707061,
    -0.00905702179254142,
    -0.021161288379106136,
    -0.014022519166323022,
    -0.004001329593238623
)

$bValues = @(
    0.22785410623173163,
    0.1414323483906026,
    0.07827534693336702,
    0.07003326479901162,
    0.06623166924227153,
    0.04079135760901664,
    0.03073799740481631,
    0.020672225126570076,
    0.018626224237523026,
    0.016623837876950986,
    0.01362110823837881,
    -0.006734929652751465,
    -0.010248701156681861,
    -0.01826002240122282,
    -0.019273655062618467,
    0.006003441757640893,
    0.003999999975077273,
    0.07993383731979264,
    0.07475435779345219,
    0.008005787088785965,
    0.003999999980056401,
    0.058442859367694666,
    0.052772868466474776,
    0.019999999931815005,
    0.023161838494262454,
    0.020648847146746974,
    0.018087312259057686,
    0.016057021759129597,
    0.009054838941513488,
    0.021022519131564277,
    0.014001329549424568,
    0.003999999973768098
)

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
